$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextCell "D2" "65.117.87"
Set-TextCell "E2" "  -2.06%  "
Set-TextCell "D3" "3.474.31"
Set-TextCell "E3" "  -1.06%  "
Set-TextCell "E4" "  +0.03%  "
Set-TextCell "D5" "587.03"
Set-TextCell "E5" "  -2.94%  "
Set-TextCell "D6" "136.70"
Set-TextCell "E6" "  -4.81%  "
Set-TextCell "D7" "3.473.44"
Set-TextCell "E7" "  -1.02%  "
Set-TextCell "E8" "  +0.08%  "
Set-TextCell "E9" "  -2.82%  "
Set-TextCell "E10" "  -5.76%  "
Set-TextCell "E11" "  -7.13%  "
Set-TextCell "E12" "  -4.77%  "
Set-TextCell "D13" "4.065.02"
Set-TextCell "E13" "  -0.88%  "
Set-TextCell "E14" "  -6.51%  "
Set-TextCell "D15" "3.477.66"
Set-TextCell "E15" "  -1.71%  "
Set-TextCell "D16" "26.56"
Set-TextCell "E16" "  -7.11%  "
Set-TextCell "E17" "  -1.35%  "
Set-TextCell "D18" "65.073.87"
Set-TextCell "E18" "  -1.92%  "
Set-TextCell "E19" "  -8.70%  "
Set-TextCell "D20" "5.76"
Set-TextCell "E20" "  -5.18%  "
Set-TextCell "E21" "  -4.52%  "
Set-TextCell "D22" "388.08"
Set-TextCell "E22" "  -7.85%  "
Set-TextCell "D23" "0.556"
Set-TextCell "E23" "  -5.15%  "
Set-TextCell "D24" "0.999"
Set-TextCell "E24" "  -0.06%  "
Set-TextCell "B25" "Litecoin"
Set-TextCell "C25" "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextCell "D25" "72.48"
Set-TextCell "E25" "  -5.48%  "
Set-TextCell "B26" "LEO"
Set-TextCell "C26" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextCell "D26" "5.75"
Set-TextCell "E26" "  -0.20%  "
Set-TextCell "D27" "3.617.23"
Set-TextCell "E27" "  -1.06%  "
Set-TextCell "E28" "  -2.10%  "
Set-TextCell "E29" "  -0.02%  "
Set-TextCell "E30" "  -5.15%  "
Set-TextCell "E31" "  -8.77%  "
Set-TextCell "E32" "  -9.80%  "
Set-TextCell "D33" "3.495.87"
Set-TextCell "E33" "  -0.61%  "
Set-TextCell "E34" "  -0.06%  "
Set-TextCell "E35" "  -6.62%  "
Set-TextCell "E36" "  -4.54%  "
Set-TextCell "D37" "170.40"
Set-TextCell "E37" "  -1.73%  "
Set-TextCell "D38" "1.19"
Set-TextCell "E38" "  -9.54%  "
Set-TextCell "D39" "6.82"
Set-TextCell "E39" "  -8.80%  "
Set-TextCell "D40" "1.47"
Set-TextCell "E40" "  -9.43%  "
Set-TextCell "E41" "  -8.68%  "
Set-TextCell "D42" "0.0776"
Set-TextCell "E42" "  -3.41%  "
Set-TextCell "E43" "  -4.56%  "
Set-TextCell "B44" "OKB"
Set-TextCell "C44" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextCell "D44" "42.56"
Set-TextCell "E44" "  -6.49%  "
Set-TextCell "B45" "FirstDigitalUSD"
Set-TextCell "C45" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextCell "D45" "1.00"
Set-TextCell "E45" "  +0.07%  "
Set-TextCell "D46" "24.96"
Set-TextCell "E46" "  +8.62%  "
Set-TextCell "E47" "  -11.98%  "
Set-TextCell "D48" "1.16"
Set-TextCell "E48" "  +4.44%  "
Set-TextCell "E49" "  -8.20%  "
Set-TextCell "E50" "  -4.58%  "
Set-TextCell "D51" "2.214.42"
Set-TextCell "E51" "  -3.64%  "
